# Fruta / hortaliza, semanal
# Insert a new weekly record as row 10, pushing the existing rows 10-85
# down to 11-86 (dimension grows from A1:R85 to A1:R86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10; this shifts every row
# from 10..85 down to 11..86 and keeps rows 1..9 untouched.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record's data.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44881
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112012
$ws.Range("G10").Value = "Espinaca"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 1400
$ws.Range("K10").Value = 1300
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = 1414
$ws.Range("N10").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 471
$ws.Range("Q10").Value = 3
$ws.Range("R10").Value = "Hortaliza"
